$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9868761301040649
$ws.Range("B1").Value = 1.902388334274292
$ws.Range("C1").Value = 5.316630363464355
$ws.Range("D1").Value = 2.256963968276978
$ws.Range("E1").Value = 1.309836506843567
